$d = $word.ActiveDocument

$replacements = @(
    @("586×7=", "615×5="),
    @("420×8=", "128×2="),
    @("998×7=", "682×4="),
    @("791×4=", "131×6="),
    @("269×3=", "774×3="),
    @("127×5=", "700×5="),
    @("636×6=", "670×4="),
    @("783×7=", "819×7="),
    @("432×9=", "900×8="),
    @("563×3=", "289×9="),
    @("665×6=", "720×8="),
    @("891×8=", "362×8="),
    @("447×2=", "698×6="),
    @("346×6=", "970×2="),
    @("759×8=", "923×5="),
    @("908×7=", "594×5="),
    @("898×7=", "621×2="),
    @("675×9=", "914×2="),
    @("923×2=", "631×3="),
    @("555×9=", "490×8="),
    @("720×7=", "959×8="),
    @("478×5=", "404×7="),
    @("840×6=", "306×8="),
    @("575×3=", "764×5="),
    @("619×9=", "179×7="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
